$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "Phones" column (E), shifting
# Phones -> G and Emails -> H. Excel will shift existing column F (Emails)
# to H automatically, and carry formatting along with it.
$ws.Range("E:F").Insert()

# Header row relabeling
$ws.Range("E1").Value = "Title"
$ws.Range("F1").Value = "Country"
$ws.Range("G1").Value = "Phone"

# Copy header style (bold/border/center) from D1 onto the two new header cells
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Title / Country per row (index by row number)
$titles = @{
    2  = "Project Director"
    3  = "Project Director"
    4  = "Project Director"
    5  = "Director of SHEQ"
    6  = "Project Director"
    7  = "Project Manager"
    8  = "Contracts Director"
    9  = "Project Director"
    10 = "Senior SHEQ Advisor"
    11 = "Construction Operations Manager"
    12 = "Project Director"
    13 = "Project Manager"
    14 = "Assistant Project Manager"
    15 = "SHEQ Officer"
    16 = "Head of Buying"
    17 = "Health and Safety Officer"
    18 = "M&E"
    19 = "Snr. SHEQ Advisor"
    20 = "Senior Construction Manager"
    21 = "Project Manager"
    22 = "Health and safety administrator"
}

$countries = @{
    2  = "United Kingdom"
    3  = "Cookstown, Northern Ireland, United Kingdom"
    4  = "Belfast"
    5  = "Cookstown"
    6  = "London, England, United Kingdom"
    7  = "London, England, United Kingdom"
    8  = "London, United Kingdom"
    9  = "Belfast Metropolitan Area"
    10 = "Newcastle upon Tyne, England, United Kingdom"
    11 = "Ireland"
    12 = "London Area, United Kingdom"
    13 = "Belfast, Northern Ireland, United Kingdom"
    14 = "London, England, United Kingdom"
    15 = "Belfast, United Kingdom"
    16 = "Cookstown, N.Ireland"
    17 = "London, United Kingdom"
    18 = "Cookstown"
    19 = "London, England, United Kingdom"
    20 = "Belfast, United Kingdom"
    21 = "Cookstown, Northern Ireland, United Kingdom"
    22 = "Wembley, England, United Kingdom"
}

foreach ($r in 2..22) {
    $ws.Cells.Item($r, 5).Value = $titles[$r]
    $ws.Cells.Item($r, 6).Value = $countries[$r]
}
